$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lama2", "Rpsa", "ECs", 3, 1, 4.608931666666666, 13.826795, 0.02269509467890621, 0.02269509467890622, 3, 1, 112.708133, 338.124399, 0.2121524692929861, 0.2121524692929861, 519.4640832745783, 4675.176749471205, 0.004814820376968063, 0.004814820376968063),
    @("ECs", "Lama2", "Rpsa", "FAPs", 3, 1, 4.608931666666666, 13.826795, 0.02269509467890621, 0.02269509467890622, 3, 1, 189.57842, 568.7352599999999, 0.3568467408440064, 0.3568467408440064, 873.753983254633, 7863.785849291698, 0.008098670569313834, 0.008098670569313836),
    @("ECs", "Lama2", "Rpsa", "M2", 3, 1, 4.608931666666666, 13.826795, 0.02269509467890621, 0.02269509467890622, 3, 1, 110.6512236666667, 331.953671, 0.2082807132576123, 0.2082807132576123, 509.983928712716, 4589.855358414445, 0.004726950507171628, 0.004726950507171628),
    @("ECs", "Lama2", "Rpsa", "sCs", 3, 1, 4.608931666666666, 13.826795, 0.02269509467890621, 0.02269509467890622, 3, 1, 118.3222806666667, 354.966842, 0.2227200766053953, 0.2227200766053952, 545.339306236821, 4908.053756131389, 0.00505465322545269, 0.00505465322545269),
    @("FAPs", "Lama2", "Rpsa", "ECs", 3, 1, 168.218394, 504.655182, 0.8283334739316415, 0.8283334739316416, 3, 1, 112.708133, 338.124399, 0.2121524692929861, 0.2121524692929861, 18959.58112399841, 170636.2301159856, 0.1757329918926351, 0.1757329918926351),
    @("FAPs", "Lama2", "Rpsa", "FAPs", 3, 1, 168.218394, 504.655182, 0.8283334739316415, 0.8283334739316416, 3, 1, 189.57842, 568.7352599999999, 0.3568467408440064, 0.3568467408440064, 31890.57734945748, 287015.1961451173, 0.2955881005045, 0.2955881005045001),
    @("FAPs", "Lama2", "Rpsa", "M2", 3, 1, 168.218394, 504.655182, 0.8283334739316415, 0.8283334739316416, 3, 1, 110.6512236666667, 331.953671, 0.2082807132576123, 0.2082807132576123, 18613.57113934146, 167522.1402540731, 0.1725258867656381, 0.1725258867656381),
    @("FAPs", "Lama2", "Rpsa", "sCs", 3, 1, 168.218394, 504.655182, 0.8283334739316415, 0.8283334739316416, 3, 1, 118.3222806666667, 354.966842, 0.2227200766053953, 0.2227200766053952, 19903.98402816392, 179135.8562534752, 0.1844864947688684, 0.1844864947688684),
    @("M2", "Lama2", "Rpsa", "ECs", 2, 0.6666666666666666, 0.1627236666666667, 0.488171, 0.00080127658394417, 0.00080127658394417, 3, 1, 112.708133, 338.124399, 0.2121524692929861, 0.2121524692929861, 18.34028066491434, 165.062525984229, 0.0001699928058704043, 0.0001699928058704043),
    @("M2", "Lama2", "Rpsa", "FAPs", 2, 0.6666666666666666, 0.1627236666666667, 0.488171, 0.00080127658394417, 0.00080127658394417, 3, 1, 189.57842, 568.7352599999999, 0.3568467408440064, 0.3568467408440064, 30.84889562327333, 277.64006060946, 0.000285932937495096, 0.000285932937495096),
    @("M2", "Lama2", "Rpsa", "M2", 2, 0.6666666666666666, 0.1627236666666667, 0.488171, 0.00080127658394417, 0.00080127658394417, 3, 1, 110.6512236666667, 331.953671, 0.2082807132576123, 0.2082807132576123, 18.00557283619345, 162.050155525741, 0.0001668904584205148, 0.0001668904584205148),
    @("M2", "Lama2", "Rpsa", "sCs", 2, 0.6666666666666666, 0.1627236666666667, 0.488171, 0.00080127658394417, 0.00080127658394417, 3, 1, 118.3222806666667, 354.966842, 0.2227200766053953, 0.2227200766053952, 19.25383535844244, 173.284518225982, 0.000178460382158155, 0.0001784603821581549),
    @("sCs", "Lama2", "Rpsa", "ECs", 3, 1, 30.09047233333333, 90.271417, 0.148170154805508, 0.148170154805508, 3, 1, 112.708133, 338.124399, 0.2121524692929861, 0.2121524692929861, 3391.440957778154, 30522.96862000339, 0.03143466421751254, 0.03143466421751254),
    @("sCs", "Lama2", "Rpsa", "FAPs", 3, 1, 30.09047233333333, 90.271417, 0.148170154805508, 0.148170154805508, 3, 1, 189.57842, 568.7352599999999, 0.3568467408440064, 0.3568467408440064, 5704.504202007046, 51340.53781806341, 0.05287403683269744, 0.05287403683269744),
    @("sCs", "Lama2", "Rpsa", "M2", 3, 1, 30.09047233333333, 90.271417, 0.148170154805508, 0.148170154805508, 3, 1, 110.6512236666667, 331.953671, 0.2082807132576123, 0.2082807132576123, 3329.547584391312, 29965.92825952181, 0.03086098552638205, 0.03086098552638204),
    @("sCs", "Lama2", "Rpsa", "sCs", 3, 1, 30.09047233333333, 90.271417, 0.148170154805508, 0.148170154805508, 3, 1, 118.3222806666667, 354.966842, 0.2227200766053953, 0.2227200766053952, 3560.373312817235, 32043.35981535511, 0.03300046822891602, 0.03300046822891602)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
